$wb = $excel.ActiveWorkbook

# Rename sheet "Energy_source" -> "Energy_source_year"
$wsEnergySource = $wb.Worksheets.Item("Energy_source")
$wsEnergySource.Name = "Energy_source_year"

# Update values on res_type_Energy_source (surface-related values divided by 2.3)
$wsRes = $wb.Worksheets.Item("res_type_Energy_source")
$wsRes.Range("D5").Value = 72.82295273457801
$wsRes.Range("D7").Value = 63.968949809466615
$wsRes.Range("D8").Value = 14.762065340646139
$wsRes.Range("D13").Value = 81.012829139343438
$wsRes.Range("D15").Value = 63.741428610609482
$wsRes.Range("D16").Value = 14.709560448602184
$wsRes.Range("D21").Value = 95.574512422360257
$wsRes.Range("D23").Value = 91.69184782608697
$wsRes.Range("D24").Value = 96.105887779722764

# Update selection on Energy_source_Vecteur sheet: F8 -> C21
$wsVecteur = $wb.Worksheets.Item("Energy_source_Vecteur")
$wsVecteur.Range("C21").Select()

# Update selection sqref on res_type_Energy_source: D2 -> D2:D25
$wsRes.Range("D2:D25").Select()

# Update selection on Energy_source_year sheet: G9 -> A27 (select last so this
# sheet ends up as the active/tabSelected sheet, matching the original file
# where "Energy_source" / "Energy_source_year" is the tabSelected sheet)
$wsEnergySource.Range("A27").Select()
